# Auto-generated edit script: updates crypto price/volume data per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D='29.788.39'; E='  -0.39%  ' },
    @{ Row=3; D='1.871.18'; E='  -0.16%  ' },
    @{ Row=4; D='0.9997'; E='  -0.09%  ' },
    @{ Row=5; D='0.7265'; E='  -1.90%  ' },
    @{ Row=6; D='241.27'; E='  -0.47%  ' },
    @{ Row=7; D='1.0000'; E='  -0.01%  ' },
    @{ Row=8; D='0.3140'; E='  -0.31%  ' },
    @{ Row=9; D='0.07152'; E='  -0.03%  ' },
    @{ Row=10; D='24.53'; E='  -0.57%  ' },
    @{ Row=11; D='0.08156'; E='  -2.88%  ' },
    @{ Row=12; B='Polygon'; C='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D='0.7448'; E='  -0.78%  ' },
    @{ Row=13; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='5.362'; E='  -1.04%  ' },
    @{ Row=14; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='1.858.61'; E='  +0.13%  ' },
    @{ Row=15; D='92.61'; E='  +0.11%  ' },
    @{ Row=16; D='29.800.39'; E='  -0.31%  ' },
    @{ Row=17; D='6.023'; E='  -1.38%  ' },
    @{ Row=18; D='247.65'; E='  +1.84%  ' },
    @{ Row=19; D='13.44'; E='  -1.04%  ' },
    @{ Row=20; D='0.000007812'; E='  -0.02%  ' },
    @{ Row=21; D='1.000'; E='  +0.15%  ' },
    @{ Row=22; D='2.127.98'; E='  +0.22%  ' },
    @{ Row=23; D='1.000'; E='  -0.09%  ' },
    @{ Row=24; D='7.749'; E='  -2.83%  ' },
    @{ Row=25; D='0.1527'; E='  -1.30%  ' },
    @{ Row=26; D='9.240'; E='  -0.54%  ' },
    @{ Row=27; D='163.81'; E='  -0.85%  ' },
    @{ Row=28; D='18.59'; E='  -0.15%  ' },
    @{ Row=29; D='2.014'; E='  -1.04%  ' },
    @{ Row=30; D='1.446'; E='  -2.92%  ' },
    @{ Row=31; D='4.538'; E='  -1.36%  ' },
    @{ Row=32; D='1.520'; E='  -0.43%  ' },
    @{ Row=33; D='4.186'; E='  -1.66%  ' },
    @{ Row=34; D='0.05409'; E='  +1.52%  ' },
    @{ Row=35; D='1.231'; E='  -0.46%  ' },
    @{ Row=36; D='0.7385'; E='  -2.18%  ' },
    @{ Row=37; E='  +0.55%  ' },
    @{ Row=38; D='2.702'; E='  +0.35%  ' },
    @{ Row=39; D='0.01930'; E='  -1.07%  ' },
    @{ Row=40; E='  -0.55%  ' },
    @{ Row=41; D='0.4479'; E='  -0.31%  ' },
    @{ Row=42; D='0.8882'; E='  +3.68%  ' },
    @{ Row=43; D='5.979'; E='  -1.25%  ' },
    @{ Row=44; D='71.37'; E='  -1.07%  ' },
    @{ Row=45; B='PaxDollar'; C='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D='1.000'; E='  -0.08%  ' },
    @{ Row=46; B='Maker'; C='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D='1.038.10'; E='  -6.60%  ' },
    @{ Row=47; D='104.00'; E='  +0.88%  ' },
    @{ Row=48; D='7.484'; E='  -2.26%  ' },
    @{ Row=49; D='1.821'; E='  -0.89%  ' },
    @{ Row=50; D='9.574'; E='  +0.55%  ' },
    @{ Row=51; D='2.030.16'; E='  +0.46%  ' },
)

$colMap = @{ B = 2; C = 3; D = 4; E = 5 }

foreach ($u in $updates) {
    $row = $u.Row
    foreach ($colLetter in $colMap.Keys) {
        if ($u.ContainsKey($colLetter)) {
            $col = $colMap[$colLetter]
            $cell = $ws.Cells.Item($row, $col)
            $cell.NumberFormat = "@"
            $cell.Value = $u[$colLetter]
            $cell.ClearFormats()
        }
    }
}
